# Add 6.9% as additional DM input (keep 5.1%)
#
# The "Slurry" sheet has a small table (A1:D5) of digestate NH3 emission
# factors broken out by acid-dosing rate, all currently computed for a
# single dry-matter (man.dm) input of 5.1%. This adds a duplicate block of
# four rows for man.dm = 6.9%, reusing the same acid-dosing labels / pH
# formulas, directly below the existing block (new rows 6-9).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slurry")

# Duplicate the existing data rows (2-5) into the new block (6-9), carrying
# over cell styles/number-formats along with values & formulas.
$ws.Range("A2:D5").Copy($ws.Range("A6"))

# Row 6 (0 kg/t acid) has no formula in column D, just the plain man.ph
# value, same as row 2 - the Copy above already handled that. Rows 7-9 use
# a formula in column D; Copy() only carried over the computed value, so
# restore the formulas explicitly (identical formulas to rows 3-5).
$ws.Range("D7").Formula = "=7.9-1.38"
$ws.Range("D8").Formula = "=7.9-0.71"
$ws.Range("D9").Formula = "=7.9-1.11"

# The only real change versus the copied rows 2-5: man.dm (column C) is
# 6.9 for the new block instead of 5.1.
$ws.Range("C6").Value = 6.9
$ws.Range("C7").Value = 6.9
$ws.Range("C8").Value = 6.9
$ws.Range("C9").Value = 6.9

# Leave the cursor where the author's saved selection shows.
$ws.Range("F13").Select() | Out-Null
